$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel columns A-F ---
$ws.Range("A1").Value = "Municipio codigo"
$ws.Range("B1").Value = "Estado civil"
$ws.Range("C1").Value = "Estado civil, código"
$ws.Range("D1").Value = "Sexo, código"
$ws.Range("E1").Value = "Sexo"
$ws.Range("F1").Value = "Municipio nombre"

# --- Data rows 2-4: values moved from C/D into E/F (swapped), C/D become "null" ---
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "dim"

$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "URI-Municipio"
